# Adapt column header formatting to respective input file names (#7)
# - rename "<field>_old" -> "<field>_FV2310" and "<field>_new" -> "<field>_FV2404"
#   for the header row (row 1)
# - turn the data range into an Excel Table ("Table1")
# - freeze the header row (split/freeze pane at row 2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header cells (A1:J1 "_old" -> "_FV2310", L1:U1 "_new" -> "_FV2404")
#    K1 ("diff") is left untouched.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# ---------------------------------------------------------------------------
# 2) Turn A1:U72 into an Excel Table.
#    ListObjects.Add derives a "header row" differential style (dxf) from
#    whatever formatting is currently applied to the header cells, which
#    would introduce a styles.xml change that isn't part of this edit.  To
#    avoid that, stash the header formatting on a scratch row, strip the
#    header formatting, create the table, then restore the original
#    formatting from the scratch row and discard the scratch row.
# ---------------------------------------------------------------------------
$header = $ws.Range("A1:U1")
$scratchRow = $ws.Rows.Item(100)
$scratch = $ws.Range("A100:U100")

$header.Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats
$header.ClearFormats()
$excel.CutCopyMode = 0

$tableRange = $ws.Range("A1:U72")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

$scratch.Copy()
$header.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$scratchRow.Delete()

# ---------------------------------------------------------------------------
# 3) Freeze the header row.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Table: $($lo.Name), Range: $($lo.Range.Address())"
